$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1767.6
$ws.Range("I6").Value = 946
$ws.Range("K6").Value = 2838
$ws.Range("M6").Value = -2726

$ws.Range("H17").Value = 1268
$ws.Range("J17").Value = 1268
$ws.Range("L17").Value = 3804
$ws.Range("N17").Value = -4140

$ws.Range("H61").Value = 7899
$ws.Range("I61").Value = 4824.75
$ws.Range("K61").Value = 14474.25
$ws.Range("M61").Value = -14302.25

$ws.Range("H97").Value = 2928
$ws.Range("J97").Value = 2928
$ws.Range("L97").Value = 8784
$ws.Range("N97").Value = -9776

$ws.Range("H99").Value = 2839.6
$ws.Range("I99").Value = 311
$ws.Range("K99").Value = 933
$ws.Range("M99").Value = 565

$ws.Range("H100").Value = 2329.0393
$ws.Range("I100").Value = 1915.4524
$ws.Range("K100").Value = 1915.4524
$ws.Range("M100").Value = -1374.4524

$ws.Range("H125").Value = 5812.125
$ws.Range("I125").Value = 3666
$ws.Range("J125").Value = 7099.8
$ws.Range("K125").Value = 32994
$ws.Range("L125").Value = 63898.2
$ws.Range("M125").Value = -30534
$ws.Range("N125").Value = -68818.20000000001

$ws.Range("H137").Value = 3267.182
$ws.Range("I137").Value = 3365.5557
$ws.Range("K137").Value = 10096.6671
$ws.Range("M137").Value = -7546.667099999999

$ws.Range("H138").Value = 3605
$ws.Range("I138").Value = 3071.6365
$ws.Range("J138").Value = 3860.087
$ws.Range("K138").Value = 9214.9095
$ws.Range("L138").Value = 11580.261
$ws.Range("M138").Value = -4074.9095
$ws.Range("N138").Value = -21860.261

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 893.34283
$ws.Range("I2").Value = 893.34283
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 893.34283
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -780.34283
$ws.Range("N2").ClearContents()

$ws.Range("H32").Value = 52956.332
$ws.Range("I32").Value = 93106
$ws.Range("K32").Value = 93106
$ws.Range("M32").Value = -92819

$ws.Range("H45").Value = 3833.5386
$ws.Range("I45").Value = 2158.8333
$ws.Range("J45").Value = 5269
$ws.Range("K45").Value = 2158.8333
$ws.Range("L45").Value = 5269
$ws.Range("M45").Value = -1781.8333
$ws.Range("N45").Value = -6023

$ws.Range("H61").Value = 2496.1177
$ws.Range("I61").Value = 2051.6667
$ws.Range("K61").Value = 2051.6667
$ws.Range("M61").Value = -1839.6667

$ws.Range("H74").Value = 2048.7
$ws.Range("I74").Value = 1775.5555
$ws.Range("K74").Value = 1775.5555
$ws.Range("M74").Value = -901.5554999999999

$ws.Range("H77").Value = 2048.7
$ws.Range("I77").Value = 1775.5555
$ws.Range("K77").Value = 8877.7775
$ws.Range("M77").Value = -4509.7775

$ws.Range("H96").Value = 27793.857
$ws.Range("I96").Value = 24993
$ws.Range("J96").Value = 28260.666
$ws.Range("K96").Value = 24993
$ws.Range("L96").Value = 28260.666
$ws.Range("M96").Value = -22247
$ws.Range("N96").Value = -33752.666

$ws.Range("H97").Value = 909.6087
$ws.Range("I97").Value = 792.9
$ws.Range("J97").Value = 999.38464
$ws.Range("K97").Value = 792.9
$ws.Range("L97").Value = 999.38464
$ws.Range("M97").Value = -296.9
$ws.Range("N97").Value = -1991.38464

$ws.Range("H105").Value = 98833.336
$ws.Range("J105").Value = 98833.336
$ws.Range("L105").Value = 98833.336
$ws.Range("N105").Value = -105821.336

$ws.Range("H116").Value = 893.34283
$ws.Range("I116").Value = 893.34283
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 893.34283
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1400.65717
$ws.Range("N116").ClearContents()

$ws.Range("H122").Value = 1746.5385
$ws.Range("I122").Value = 1256.4445
$ws.Range("K122").Value = 3769.3335
$ws.Range("M122").Value = -1319.3335

$ws.Range("H136").Value = 2496.1177
$ws.Range("I136").Value = 2051.6667
$ws.Range("K136").Value = 6155.000100000001
$ws.Range("M136").Value = -3605.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 893.34283
$ws.Range("I3").Value = 893.34283
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 893.34283
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -779.34283
$ws.Range("N3").ClearContents()

$ws.Range("H94").Value = 1528.8667
$ws.Range("I94").Value = 1078.75
$ws.Range("K94").Value = 1078.75
$ws.Range("M94").Value = -627.75

$ws.Range("H95").Value = 18373.4
$ws.Range("J95").Value = 18373.4
$ws.Range("L95").Value = 18373.4
$ws.Range("N95").Value = -23865.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 5061.3335
$ws.Range("J96").Value = 5061.3335
$ws.Range("L96").Value = 5061.3335
$ws.Range("N96").Value = -10553.3335

$ws.Range("H138").Value = 79999
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 79999
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 79999
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -90279

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 5352.75
$ws.Range("J64").Value = 6299.6665
$ws.Range("L64").Value = 18898.9995
$ws.Range("N64").Value = -19438.9995

$ws.Range("H67").Value = 5352.75
$ws.Range("J67").Value = 6299.6665
$ws.Range("L67").Value = 18898.9995
$ws.Range("N67").Value = -20770.9995

$ws.Range("H97").Value = 448.27274
$ws.Range("I97").Value = 285.25
$ws.Range("J97").Value = 541.4286
$ws.Range("K97").Value = 855.75
$ws.Range("L97").Value = 1624.2858
$ws.Range("M97").Value = -359.75
$ws.Range("N97").Value = -2616.2858

$ws.Range("H132").Value = 3357.5715
$ws.Range("I132").Value = 3474.5
$ws.Range("J132").Value = 3201.6667
$ws.Range("K132").Value = 31270.5
$ws.Range("L132").Value = 28815.0003
$ws.Range("M132").Value = -28740.5
$ws.Range("N132").Value = -33875.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11047.55
$ws.Range("I80").Value = 7730.1333
$ws.Range("K80").Value = 7730.1333
$ws.Range("M80").Value = -6732.1333

$ws.Range("H83").Value = 11047.55
$ws.Range("I83").Value = 7730.1333
$ws.Range("K83").Value = 38650.6665
$ws.Range("M83").Value = -33658.6665

$ws.Range("H92").Value = 18588.928
$ws.Range("I92").Value = 5990
$ws.Range("J92").Value = 19558.076
$ws.Range("K92").Value = 5990
$ws.Range("L92").Value = 19558.076
$ws.Range("M92").Value = -4118
$ws.Range("N92").Value = -23302.076

$ws.Range("H107").Value = 50321.24
$ws.Range("J107").Value = 2064.8
$ws.Range("L107").Value = 2064.8
$ws.Range("N107").Value = -5904.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2201.3684
$ws.Range("I82").Value = 783
$ws.Range("K82").Value = 783
$ws.Range("M82").Value = -422

$ws.Range("H85").Value = 2201.3684
$ws.Range("I85").Value = 783
$ws.Range("K85").Value = 783
$ws.Range("M85").Value = 465

$ws.Range("H93").Value = 2055.0625
$ws.Range("I93").Value = 2074.4285
$ws.Range("K93").Value = 2074.4285
$ws.Range("M93").Value = -826.4285

$ws.Range("H122").Value = 11250
$ws.Range("J122").Value = 12500
$ws.Range("L122").Value = 37500
$ws.Range("N122").Value = -42400

$ws.Range("H132").Value = 51853.125
$ws.Range("I132").Value = 53870.914
$ws.Range("K132").Value = 161612.742
$ws.Range("M132").Value = -159082.742

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H81").Value = 2522.7778
$ws.Range("I81").Value = 2873.5715
$ws.Range("K81").Value = 5747.143
$ws.Range("M81").Value = -4686.143

$ws.Range("H84").Value = 2522.7778
$ws.Range("I84").Value = 2873.5715
$ws.Range("K84").Value = 28735.715
$ws.Range("M84").Value = -23431.715

$ws.Range("H100").Value = 860.2632
$ws.Range("I100").Value = 764.2857
$ws.Range("K100").Value = 1528.5714
$ws.Range("M100").Value = -987.5714

$ws.Range("H122").Value = 4513.3335
$ws.Range("J122").Value = 3750
$ws.Range("L122").Value = 11250
$ws.Range("N122").Value = -16150
